# Fix broken / mismatched GitHub links in the portfolio database sheet.
#
# Rows (by project, in their current sheet order):
#   Row 2 - Movie Recommender System with Sentiment Analysis
#   Row 3 - Portfolio Website                -> F3 had the WRONG url (Airbnb repo) and no hyperlink
#   Row 4 - Human Rights First Asylum         -> F4 pointed at the old Lambda-School-Labs fork, no hyperlink
#   Row 5 - Airbnb in Los Angeles             -> F5 had the WRONG url (Movie Recommender repo)
#   Row 6 - Image Classifier using VGG-19 CNN
#   Row 7 - Airbus / Data Storytelling        -> F7 had the WRONG url (human-rights-first repo)
#   Row 8 - Medical Cannabis Recommender System
#   Row 9 - Kickstarter Success Classifier

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Airbnb in Los Angeles -> correct github link, update existing hyperlink target ---
$ws.Range("F5").Value = "https://github.com/navroz-lamba/Prediciting-Airbnb-prices-in-LA"
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq "`$F`$5") {
        $h.Address = "https://github.com/navroz-lamba/Prediciting-Airbnb-prices-in-LA"
    }
}

# --- Row 4: Human Rights First Asylum -> correct github link (new owner), add missing hyperlink ---
$ws.Range("F4").Value = "https://github.com/navroz-lamba/human-rights-first-asylum-ds-a"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/navroz-lamba/human-rights-first-asylum-ds-a")
$ws.Range("F4").Style = "Hyperlink"

# --- Row 7: Airbus, the new King of the Skies?! | Data Storytelling -> correct github link ---
$ws.Range("F7").Value = "https://github.com/navroz-lamba/DataStorytelling"
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq "`$F`$7") {
        $h.Address = "https://github.com/navroz-lamba/DataStorytelling"
    }
}

# --- Row 3: Portfolio Website -> correct github link, add missing hyperlink ---
$ws.Range("F3").Value = "https://github.com/navroz-lamba/Portfolio-Website"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/navroz-lamba/Portfolio-Website")
$ws.Range("F3").Style = "Hyperlink"

# Leave the cursor on the cell that was actually fixed last.
$ws.Range("F5").Select()

$wb.Save()
